$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 122, shifting existing rows 122:193 down to 123:194
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row 122 with the new data record
$ws.Range("A122").Value = 1
$ws.Range("B122").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C122").Value = "Arica y Parinacota"
$ws.Range("D122").Value = 45233
$ws.Range("E122").Value = 15
$ws.Range("F122").Value = "Fruta"
$ws.Range("G122").Value = 100102
$ws.Range("H122").Value = "Cítricos"
$ws.Range("I122").Value = 100102005
$ws.Range("J122").Value = "Naranja"
$ws.Range("K122").Value = "Lane Late"
$ws.Range("L122").Value = "Segunda"
$ws.Range("M122").Value = 400
$ws.Range("N122").Value = 800
$ws.Range("O122").Value = 850
$ws.Range("P122").Value = 819
$ws.Range("Q122").Value = "`$/kilo (en caja de 20 kilos)"
$ws.Range("R122").Value = "Región de O'Higgins"
$ws.Range("S122").Value = 819
$ws.Range("T122").Value = 1
